$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
Write-Host "Panes.Count" $excel.ActiveWindow.Panes.Count
for ($i=1; $i -le $excel.ActiveWindow.Panes.Count; $i++) {
    $p = $excel.ActiveWindow.Panes.Item($i)
    Write-Host "pane $i scrollcol" $p.ScrollColumn "scrollrow" $p.ScrollRow
}
